# Update existing rows 2-4 and append new rows 5-8 in the 'Probabilities' sheet
# reflecting the latest KHL tour fixtures/odds/probabilities refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Probabilities")

$rowsData = @{}
$rowsData[2] = @(1369, "2025-12-05T17:00:00", "Трактор", "ЦСКА", 897850, "https://text.khl.ru/text/897850.html", 5.71875, 1.312362, 3.034944, 1.4375, 3.578125, 2.173653, 7.031112, 36.896451, 23.317432, 60.213883, 0.2, -0.160025, 0.64388, 0.145519, 0.209371, 0.174766, 0.824004, 0.319654, 0.679116, 0.486326, 0.512444, 0.646103, 0.352667, 0.77739, 0.221381, 0.872142, 0.127858, 0.693361, 0.306639, 0.638961, 0.361039, 0.370213, 0.629787, 0.892073, 0.523904)
$rowsData[3] = @(1369, "2025-12-05T17:00:00", "Салават Юлаев", "Ак Барс", 897851, "https://text.khl.ru/text/897851.html", 3.104859, 1.541765, 2.913973, 1.066667, 2.085763, 2.227869, 4.646624, 27.108074, 26.186936, 53.29501, -0.017864, -0.168735, 0.374402, 0.198276, 0.427287, 0.374708, 0.625258, 0.567805, 0.432161, 0.734394, 0.265571, 0.854162, 0.145804, 0.927967, 0.071999, 0.61671, 0.38329, 0.346523, 0.653477, 0.652172, 0.347828, 0.384749, 0.615251, 0.751882, 0.793337)
$rowsData[4] = @(1369, "2025-12-05T17:00:00", "Автомобилист", "Северсталь", 897855, "https://text.khl.ru/text/897855.html", 5.25, 1.59375, 3.75, 1.1875, 3.21875, 2.671875, 6.84375, 32.058457, 24.168153, 56.226611, 0.2, -0.2, 0.50456, 0.164426, 0.330383, 0.161234, 0.838135, 0.299963, 0.699406, 0.463404, 0.535965, 0.623865, 0.375504, 0.758896, 0.240473, 0.831229, 0.168771, 0.623996, 0.376004, 0.746191, 0.253809, 0.499461, 0.500539, 0.805343, 0.659074)
$rowsData[5] = @(1369, "2025-12-05T17:30:00", "Барыс", "Спартак", 897856, "https://text.khl.ru/text/897856.html", 1.394737, 4.961538, 4.421053, 5.884615, 3.639676, 4.691296, 6.356275, 23.887111, 38.240211, 62.127322, -0.2, 0.2, 0.29229, 0.131939, 0.56549, 0.033828, 0.955892, 0.082187, 0.907533, 0.162762, 0.826957, 0.274641, 0.715079, 0.407792, 0.581928, 0.878158, 0.121842, 0.704216, 0.295784, 0.947784, 0.052216, 0.846823, 0.153177, 0.564631, 0.806358)
$rowsData[6] = @(1369, "2025-12-05T18:00:00", "Лада", "Динамо М", 897854, "https://text.khl.ru/text/897854.html", 0.923077, 2.631811, 1.795369, 1.1875, 1.055288, 2.21359, 3.554888, 21.218038, 31.26235, 52.480388, -0.2, -0.056253, 0.169416, 0.194473, 0.63609, 0.587223, 0.412756, 0.768243, 0.231736, 0.88659, 0.113389, 0.951067, 0.048912, 0.981177, 0.018802, 0.284571, 0.715429, 0.090747, 0.909253, 0.64873, 0.35127, 0.380928, 0.619072, 0.593186, 0.939876)
$rowsData[7] = @(1369, "2025-12-05T19:00:00", "Нефтехимик", "Металлург Мг", 897853, "https://text.khl.ru/text/897853.html", 4.05, 5.71875, 4.875, 3.65625, 3.853125, 5.296875, 9.76875, 34.193418, 41.360427, 75.553845, 0.2, 0.2, 0.256289, 0.120081, 0.601605, 0.019086, 0.958888, 0.050109, 0.927865, 0.10688, 0.871094, 0.193457, 0.784518, 0.306625, 0.67135, 0.897049, 0.102951, 0.739576, 0.260424, 0.96847, 0.03153, 0.898227, 0.101773, 0.50909, 0.818231)
$rowsData[8] = @(1369, "2025-12-05T19:30:00", "Драконы", "Торпедо", 897852, "https://text.khl.ru/text/897852.html", 2.330401, 1.382353, 4.038462, 1.446421, 1.888411, 2.710407, 3.712754, 26.821065, 26.902948, 53.724013, -0.081103, -0.2, 0.264389, 0.179804, 0.555678, 0.325899, 0.673972, 0.513456, 0.486415, 0.685963, 0.313908, 0.818185, 0.181686, 0.905052, 0.094819, 0.562949, 0.437051, 0.293152, 0.706848, 0.753222, 0.246778, 0.508921, 0.491079, 0.633982, 0.867713)

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}

